$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 208.296368342823
$ws.Range("C5").Value = 219.166666666667
$ws.Range("I5").Value = 199.633333333333
$ws.Range("H6").Value = 4.64972181101283
